$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("1er Parcial")

$sheet1.Range("F251").Value = 7
$sheet1.Range("H251").Value = 30.4
$sheet1.Range("I251").Value = 7.8
$sheet1.Range("J251").Value = 0
$sheet1.Range("K251").Value = 0
$sheet1.Range("F253").Value = 7
$sheet1.Range("H253").Value = 17.5
$sheet1.Range("I253").Value = 8.6
$sheet1.Range("J253").Value = 0
$sheet1.Range("K253").Value = 0

$sheet2 = $wb.Worksheets.Item("2o Parcial")

$sheet2.Range("E8").Value = 22
$sheet2.Range("F8").Value = 6
$sheet2.Range("G8").Value = 78.59999999999999
$sheet2.Range("H8").Value = 21.4
$sheet2.Range("I8").Value = 7.5
$sheet2.Range("I8").NumberFormat = "0.0"
$sheet2.Range("J8").Value = 0
$sheet2.Range("K8").Value = 0
$sheet2.Range("E13").Value = 117
$sheet2.Range("F13").Value = 30
$sheet2.Range("G13").Value = 79.59999999999999
$sheet2.Range("H13").Value = 20.4
$sheet2.Range("J13").Value = 0
$sheet2.Range("K13").Value = 0
$sheet2.Range("E127").Value = 14
$sheet2.Range("F127").Value = 1
$sheet2.Range("G127").Value = 93.3
$sheet2.Range("H127").Value = 6.7
$sheet2.Range("I127").Value = 8.300000000000001
$sheet2.Range("I127").NumberFormat = "0.0"
$sheet2.Range("J127").Value = 0
$sheet2.Range("K127").Value = 0
$sheet2.Range("E128").Value = 24
$sheet2.Range("F128").Value = 1
$sheet2.Range("G128").Value = 96
$sheet2.Range("H128").Value = 4
$sheet2.Range("I128").Value = 8.4
$sheet2.Range("I128").NumberFormat = "0.0"
$sheet2.Range("J128").Value = 0
$sheet2.Range("K128").Value = 0
$sheet2.Range("E129").Value = 38
$sheet2.Range("F129").Value = 2
$sheet2.Range("G129").Value = 95
$sheet2.Range("H129").Value = 5
$sheet2.Range("I129").Value = 8.4
$sheet2.Range("I129").NumberFormat = "0.0"
$sheet2.Range("J129").Value = 0
$sheet2.Range("K129").Value = 0
$sheet2.Range("E160").Value = 34
$sheet2.Range("F160").Value = 1
$sheet2.Range("G160").Value = 97.09999999999999
$sheet2.Range("H160").Value = 2.9
$sheet2.Range("I160").Value = 6.9
$sheet2.Range("I160").NumberFormat = "0.0"
$sheet2.Range("J160").Value = 0
$sheet2.Range("K160").Value = 0
$sheet2.Range("E161").Value = 34
$sheet2.Range("F161").Value = 1
$sheet2.Range("G161").Value = 97.09999999999999
$sheet2.Range("H161").Value = 2.9
$sheet2.Range("I161").Value = 6.9
$sheet2.Range("I161").NumberFormat = "0.0"
$sheet2.Range("J161").Value = 0
$sheet2.Range("K161").Value = 0
$sheet2.Range("E174").Value = 20
$sheet2.Range("F174").Value = 4
$sheet2.Range("G174").Value = 83.3
$sheet2.Range("H174").Value = 16.7
$sheet2.Range("I174").Value = 7.6
$sheet2.Range("I174").NumberFormat = "0.0"
$sheet2.Range("J174").Value = 0
$sheet2.Range("K174").Value = 0
$sheet2.Range("J180").Value = 0
$sheet2.Range("K180").Value = 0
$sheet2.Range("E181").Value = 182
$sheet2.Range("F181").Value = 9
$sheet2.Range("G181").Value = 95.3
$sheet2.Range("H181").Value = 4.7
$sheet2.Range("I181").Value = 8.699999999999999
$sheet2.Range("J181").Value = 0
$sheet2.Range("K181").Value = 0
$sheet2.Range("E251").Value = 20
$sheet2.Range("F251").Value = 3
$sheet2.Range("G251").Value = 87
$sheet2.Range("H251").Value = 13
$sheet2.Range("I251").Value = 8.300000000000001
$sheet2.Range("I251").NumberFormat = "0.0"
$sheet2.Range("J251").Value = 0
$sheet2.Range("K251").Value = 0
$sheet2.Range("E253").Value = 37
$sheet2.Range("F253").Value = 3
$sheet2.Range("G253").Value = 92.5
$sheet2.Range("H253").Value = 7.5
$sheet2.Range("I253").Value = 8.699999999999999
$sheet2.Range("J253").Value = 0
$sheet2.Range("K253").Value = 0
$sheet2.Range("E259").Value = 11
$sheet2.Range("F259").Value = 0
$sheet2.Range("G259").Value = 100
$sheet2.Range("H259").Value = 0
$sheet2.Range("I259").Value = 8.6
$sheet2.Range("I259").NumberFormat = "0.0"
$sheet2.Range("J259").Value = 0
$sheet2.Range("K259").Value = 0
$sheet2.Range("E260").Value = 16
$sheet2.Range("F260").Value = 1
$sheet2.Range("G260").Value = 94.09999999999999
$sheet2.Range("H260").Value = 5.9
$sheet2.Range("I260").Value = 8.5
$sheet2.Range("I260").NumberFormat = "0.0"
$sheet2.Range("J260").Value = 0
$sheet2.Range("K260").Value = 0
$sheet2.Range("E261").Value = 24
$sheet2.Range("F261").Value = 1
$sheet2.Range("G261").Value = 96
$sheet2.Range("H261").Value = 4
$sheet2.Range("I261").Value = 9.4
$sheet2.Range("I261").NumberFormat = "0.0"
$sheet2.Range("J261").Value = 0
$sheet2.Range("K261").Value = 0
$sheet2.Range("E262").Value = 51
$sheet2.Range("F262").Value = 2
$sheet2.Range("G262").Value = 96.2
$sheet2.Range("H262").Value = 3.8
$sheet2.Range("I262").Value = 8.800000000000001
$sheet2.Range("I262").NumberFormat = "0.0"
$sheet2.Range("J262").Value = 0
$sheet2.Range("K262").Value = 0
$sheet2.Range("E268").Value = 39
$sheet2.Range("F268").Value = 0
$sheet2.Range("G268").Value = 100
$sheet2.Range("H268").Value = 0
$sheet2.Range("I268").Value = 9.300000000000001
$sheet2.Range("J268").Value = 0
$sheet2.Range("K268").Value = 0
$sheet2.Range("E269").Value = 39
$sheet2.Range("F269").Value = 0
$sheet2.Range("G269").Value = 100
$sheet2.Range("H269").Value = 0
$sheet2.Range("I269").Value = 9.300000000000001
$sheet2.Range("J269").Value = 0
$sheet2.Range("K269").Value = 0
$sheet2.Range("E284").Value = 18
$sheet2.Range("F284").Value = 6
$sheet2.Range("G284").Value = 75
$sheet2.Range("H284").Value = 25
$sheet2.Range("I284").Value = 7.8
$sheet2.Range("I284").NumberFormat = "0.0"
$sheet2.Range("J284").Value = 0
$sheet2.Range("K284").Value = 0
$sheet2.Range("E285").Value = 32
$sheet2.Range("F285").Value = 2
$sheet2.Range("G285").Value = 94.09999999999999
$sheet2.Range("H285").Value = 5.9
$sheet2.Range("I285").Value = 8.5
$sheet2.Range("I285").NumberFormat = "0.0"
$sheet2.Range("J285").Value = 0
$sheet2.Range("K285").Value = 0
$sheet2.Range("E286").Value = 31
$sheet2.Range("F286").Value = 0
$sheet2.Range("G286").Value = 100
$sheet2.Range("H286").Value = 0
$sheet2.Range("I286").Value = 9
$sheet2.Range("I286").NumberFormat = "0.0"
$sheet2.Range("J286").Value = 0
$sheet2.Range("K286").Value = 0
$sheet2.Range("E287").Value = 40
$sheet2.Range("F287").Value = 0
$sheet2.Range("G287").Value = 100
$sheet2.Range("H287").Value = 0
$sheet2.Range("I287").Value = 9.300000000000001
$sheet2.Range("I287").NumberFormat = "0.0"
$sheet2.Range("J287").Value = 0
$sheet2.Range("K287").Value = 0
$sheet2.Range("E288").Value = 19
$sheet2.Range("F288").Value = 5
$sheet2.Range("G288").Value = 79.2
$sheet2.Range("H288").Value = 20.8
$sheet2.Range("I288").Value = 7.6
$sheet2.Range("I288").NumberFormat = "0.0"
$sheet2.Range("J288").Value = 0
$sheet2.Range("K288").Value = 0
$sheet2.Range("E289").Value = 35
$sheet2.Range("F289").Value = 3
$sheet2.Range("G289").Value = 92.09999999999999
$sheet2.Range("H289").Value = 7.9
$sheet2.Range("I289").Value = 8.9
$sheet2.Range("I289").NumberFormat = "0.0"
$sheet2.Range("J289").Value = 0
$sheet2.Range("K289").Value = 0
$sheet2.Range("E290").Value = 175
$sheet2.Range("F290").Value = 16
$sheet2.Range("G290").Value = 91.59999999999999
$sheet2.Range("H290").Value = 8.4
$sheet2.Range("I290").Value = 8.5
$sheet2.Range("I290").NumberFormat = "0.0"
$sheet2.Range("J290").Value = 0
$sheet2.Range("K290").Value = 0
$sheet2.Range("E305").Value = 18
$sheet2.Range("F305").Value = 18
$sheet2.Range("G305").Value = 50
$sheet2.Range("H305").Value = 50
$sheet2.Range("I305").Value = 7
$sheet2.Range("I305").NumberFormat = "0.0"
$sheet2.Range("J305").Value = 0
$sheet2.Range("K305").Value = 0
$sheet2.Range("E306").Value = 20
$sheet2.Range("F306").Value = 8
$sheet2.Range("G306").Value = 71.40000000000001
$sheet2.Range("H306").Value = 28.6
$sheet2.Range("I306").Value = 7.3
$sheet2.Range("I306").NumberFormat = "0.0"
$sheet2.Range("J306").Value = 0
$sheet2.Range("K306").Value = 0
$sheet2.Range("E307").Value = 12
$sheet2.Range("F307").Value = 11
$sheet2.Range("G307").Value = 52.2
$sheet2.Range("H307").Value = 47.8
$sheet2.Range("I307").Value = 6.9
$sheet2.Range("I307").NumberFormat = "0.0"
$sheet2.Range("J307").Value = 0
$sheet2.Range("K307").Value = 0
$sheet2.Range("E308").Value = 20
$sheet2.Range("F308").Value = 10
$sheet2.Range("G308").Value = 66.7
$sheet2.Range("H308").Value = 33.3
$sheet2.Range("I308").Value = 7.2
$sheet2.Range("I308").NumberFormat = "0.0"
$sheet2.Range("J308").Value = 0
$sheet2.Range("K308").Value = 0
$sheet2.Range("E309").Value = 7
$sheet2.Range("F309").Value = 4
$sheet2.Range("G309").Value = 63.6
$sheet2.Range("H309").Value = 36.4
$sheet2.Range("I309").Value = 6.5
$sheet2.Range("I309").NumberFormat = "0.0"
$sheet2.Range("J309").Value = 0
$sheet2.Range("K309").Value = 0
$sheet2.Range("E310").Value = 77
$sheet2.Range("F310").Value = 51
$sheet2.Range("G310").Value = 60.2
$sheet2.Range("H310").Value = 39.8
$sheet2.Range("I310").Value = 7
$sheet2.Range("I310").NumberFormat = "0.0"
$sheet2.Range("J310").Value = 0
$sheet2.Range("K310").Value = 0

$sheet3 = $wb.Worksheets.Item("Final")

$sheet3.Range("E8").Value = 22
$sheet3.Range("F8").Value = 6
$sheet3.Range("G8").Value = 78.59999999999999
$sheet3.Range("H8").Value = 21.4
$sheet3.Range("I8").Value = 6.8
$sheet3.Range("E13").Value = 117
$sheet3.Range("F13").Value = 30
$sheet3.Range("G13").Value = 79.59999999999999
$sheet3.Range("H13").Value = 20.4
$sheet3.Range("I13").Value = 6.9
$sheet3.Range("I127").Value = 8.199999999999999
$sheet3.Range("E128").Value = 24
$sheet3.Range("F128").Value = 1
$sheet3.Range("G128").Value = 96
$sheet3.Range("H128").Value = 4
$sheet3.Range("I128").Value = 8.4
$sheet3.Range("E129").Value = 38
$sheet3.Range("F129").Value = 2
$sheet3.Range("G129").Value = 95
$sheet3.Range("H129").Value = 5
$sheet3.Range("I129").Value = 8.300000000000001
$sheet3.Range("I160").Value = 7.3
$sheet3.Range("I161").Value = 7.3
$sheet3.Range("E251").Value = 20
$sheet3.Range("F251").Value = 3
$sheet3.Range("G251").Value = 87
$sheet3.Range("H251").Value = 13
$sheet3.Range("I251").Value = 8.1
$sheet3.Range("E253").Value = 37
$sheet3.Range("F253").Value = 3
$sheet3.Range("G253").Value = 92.5
$sheet3.Range("H253").Value = 7.5
$sheet3.Range("I253").Value = 8.800000000000001
$sheet3.Range("I259").Value = 8.6
$sheet3.Range("I260").Value = 8.1
$sheet3.Range("I261").Value = 8.5
$sheet3.Range("I262").Value = 8.4
$sheet3.Range("E268").Value = 39
$sheet3.Range("F268").Value = 0
$sheet3.Range("G268").Value = 100
$sheet3.Range("H268").Value = 0
$sheet3.Range("I268").Value = 9.300000000000001
$sheet3.Range("E269").Value = 39
$sheet3.Range("F269").Value = 0
$sheet3.Range("G269").Value = 100
$sheet3.Range("H269").Value = 0
$sheet3.Range("I269").Value = 9.300000000000001
$sheet3.Range("I284").Value = 7.3
$sheet3.Range("I285").Value = 8.4
$sheet3.Range("E286").Value = 31
$sheet3.Range("F286").Value = 0
$sheet3.Range("G286").Value = 100
$sheet3.Range("H286").Value = 0
$sheet3.Range("I286").Value = 9
$sheet3.Range("I287").Value = 9.300000000000001
$sheet3.Range("E288").Value = 19
$sheet3.Range("F288").Value = 5
$sheet3.Range("G288").Value = 79.2
$sheet3.Range("H288").Value = 20.8
$sheet3.Range("I288").Value = 7.5
$sheet3.Range("E289").Value = 35
$sheet3.Range("F289").Value = 3
$sheet3.Range("G289").Value = 92.09999999999999
$sheet3.Range("H289").Value = 7.9
$sheet3.Range("I289").Value = 9
$sheet3.Range("E290").Value = 175
$sheet3.Range("F290").Value = 16
$sheet3.Range("G290").Value = 91.59999999999999
$sheet3.Range("H290").Value = 8.4
$sheet3.Range("I290").Value = 8.4
$sheet3.Range("E305").Value = 18
$sheet3.Range("F305").Value = 18
$sheet3.Range("G305").Value = 50
$sheet3.Range("H305").Value = 50
$sheet3.Range("I305").Value = 6.2
$sheet3.Range("E306").Value = 20
$sheet3.Range("F306").Value = 8
$sheet3.Range("G306").Value = 71.40000000000001
$sheet3.Range("H306").Value = 28.6
$sheet3.Range("I306").Value = 7
$sheet3.Range("E307").Value = 12
$sheet3.Range("F307").Value = 11
$sheet3.Range("G307").Value = 52.2
$sheet3.Range("H307").Value = 47.8
$sheet3.Range("I307").Value = 6.4
$sheet3.Range("E308").Value = 20
$sheet3.Range("F308").Value = 10
$sheet3.Range("G308").Value = 66.7
$sheet3.Range("H308").Value = 33.3
$sheet3.Range("I308").Value = 6.9
$sheet3.Range("I309").Value = 6.5
$sheet3.Range("E310").Value = 77
$sheet3.Range("F310").Value = 51
$sheet3.Range("G310").Value = 60.2
$sheet3.Range("H310").Value = 39.8
$sheet3.Range("I310").Value = 6.6
